# "feat: poprawa paneli lekarz oraz recepcja home"
#
# Piotr Bistyga's work-log block (columns K/L/M, under the "Piotr Bistyga"
# header) gets six new entries dated 2025-09-11 (serial 45911), recording
# work on the "lekarz-home" and "recepcja-home" Angular component panels
# (ts/html/css files) plus touch-ups to login.component.ts and
# home.component.css.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K65:K70 — new "Data" entries (2025-09-11), reuse the existing
#     short-date style (same style already used by the rest of column K,
#     e.g. K64) instead of letting Excel mint a brand-new number format.
$ws.Range("K64").Copy() | Out-Null
$ws.Range("K65:K70").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("K65").Value = 45911
$ws.Range("K66").Value = 45911
$ws.Range("K67").Value = 45911
$ws.Range("K68").Value = 45911
$ws.Range("K69").Value = 45911
$ws.Range("K70").Value = 45911

# --- L65:L69 — "Plik" file names, plain (default) style.
$ws.Range("L65").Value = "lekarz-home.component.html"
$ws.Range("L66").Value = "lekarz-home.component.ts"
$ws.Range("L67").Value = "recepcja-home.component.html"
$ws.Range("L68").Value = "recepcja-home.component.ts"
$ws.Range("L69").Value = "login.component.ts"

# --- L70 — same row also used the smaller "Segoe UI" style (s=7), matching
#     the style already used for F67:F70 in this block; copy it across.
$ws.Range("F67").Copy() | Out-Null
$ws.Range("L70").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$ws.Range("L70").Value = "home.component.css"

# --- M65:M70 — "Linie" line counts for each entry above.
$ws.Range("M65").Value = 10
$ws.Range("M66").Value = 6
$ws.Range("M67").Value = 10
$ws.Range("M68").Value = 15
$ws.Range("M69").Value = 3
$ws.Range("M70").Value = 7

# Column L ("Plik") now holds noticeably longer file names
# (recepcja-home.component.ts etc.) — widen it to fit, as Excel/the author
# did on save.
$ws.Columns.Item(12).ColumnWidth = 24.5

# Restore the author's last selection/cursor position.
$ws.Range("L71").Select() | Out-Null
